$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "Cleidson-Circuitos elétri"

$ws.Range("B3").Value = "Cleidson-Circuitos elétri"
$ws.Range("C3").Value = "[-, -, -, Valmir-Metrologia]"
$ws.Range("D3").Value = "-"
$ws.Range("F3").Value = "-"

$ws.Range("B4").Value = "Anselmo-Gestão integr"
$ws.Range("C4").Value = "[-, -, -, Valmir-Metrologia]"
$ws.Range("D4").Value = "-"
$ws.Range("E4").Value = "[Joel L.-Tecnologia dos Materiais., André Guimarães-Desenho Técnico]"
$ws.Range("F4").Value = "-"

$ws.Range("B6").Value = "Anselmo-Gestão integr"
$ws.Range("C6").Value = "[-, -, -, Valmir-Metrologia]"
$ws.Range("E6").Value = "[Joel L.-Tecnologia dos Materiais., -]"
$ws.Range("F6").Value = "-"

$ws.Range("C7").Value = "[-, -, -, Valmir-Metrologia]"
$ws.Range("E7").Value = "-"
$ws.Range("F7").Value = "-"

$ws.Range("B8").Value = "Cleidson-Circuitos elétri"
